# Refresh "particelle non trovate" data: the rows that used to sit at
# positions 9-18, 41 and 61-63 of the old export are gone from the refreshed
# source, so every surviving record's codice_particella (col B) /
# codice_comune_catastale (col C) slides up to fill the gap, while the
# running index kept in column A stays put on its original row. The trailing
# rows that no longer carry any data (58-71) are removed outright, shrinking
# the sheet from A1:C71 down to A1:C57.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowNums = @(9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 54, 55, 56, 57)

$bValues = @('1117/2', '1230/100', '1230/115', '1230/85', '1230/86', '1230/87', '1230/88', '1303/1', '1303/2', '1309', '1330', '1334', '1346', '1369/1', '194/4', '254/2', '337/5', '393/1', '393/2', '393/3', '465', '614', '420/101', '420/102', '420/106', '420/107', '420/109', '420/110', '420/80', '420/92', '420/93', '420/94', '420/95', '420/96', '420/97', '454', '53', '705/11', '756', '798/3', '4523', '3597/16', '3597/18', '3597/22', '2331/38', '2068/43', '2822/12', '2822/16', '2020/14')

$cValues = @(193, 193, 193, 193, 193, 193, 193, 193, 193, 193, 193, 193, 193, 193, 193, 193, 193, 193, 193, 193, 193, 193, 215, 215, 215, 215, 215, 215, 215, 215, 215, 215, 215, 215, 215, 215, 215, 215, 215, 215, 404, 9, 9, 9, 258, 310, 310, 310, 310)

for ($i = 0; $i -lt $rowNums.Count; $i++) {
    $r = $rowNums[$i]
    $ws.Cells.Item($r, 2).Value = $bValues[$i]
    $ws.Cells.Item($r, 3).Value = $cValues[$i]
}

# Rows 58-71 no longer correspond to any particella in the refreshed export;
# remove them entirely so the sheet shrinks to A1:C57.
$ws.Range("58:71").Delete()
